# Update the Delaware_B team-specific transition matrix on Sheet1 with
# refreshed probabilities (added more games, sped up simulate game logic,
# and drafted optimization logic).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2095238095238095
$ws.Range("C2").Value = 0.5180952380952381
$ws.Range("J2").Value = 0.01523809523809524
$ws.Range("P2").Value = 0.1523809523809524
$ws.Range("S2").Value = 0.1047619047619048
$ws.Range("B3").Value = 0.01805054151624549
$ws.Range("C3").Value = 0.01444043321299639
$ws.Range("J3").Value = 0.02527075812274368
$ws.Range("P3").Value = 0.7545126353790613
$ws.Range("S3").Value = 0.1877256317689531
$ws.Range("J4").Value = 0.07575757575757576
$ws.Range("P4").Value = 0.5909090909090909
$ws.Range("S4").Value = 0.3333333333333333
$ws.Range("B6").Value = 0.05882352941176471
$ws.Range("D6").Value = 0.009411764705882352
$ws.Range("F6").Value = 0.05882352941176471
$ws.Range("J6").Value = 0.2611764705882353
$ws.Range("O6").Value = 0.02352941176470588
$ws.Range("Q6").Value = 0.1458823529411765
$ws.Range("R6").Value = 0.06352941176470588
$ws.Range("S6").Value = 0.3788235294117647
$ws.Range("B7").Value = 0.1082474226804124
$ws.Range("D7").Value = 0.007731958762886598
$ws.Range("E7").Value = 0.002577319587628866
$ws.Range("F7").Value = 0.05154639175257732
$ws.Range("J7").Value = 0.1134020618556701
$ws.Range("O7").Value = 0.03092783505154639
$ws.Range("Q7").Value = 0.1907216494845361
$ws.Range("R7").Value = 0.06958762886597938
$ws.Range("S7").Value = 0.4252577319587629
$ws.Range("B8").Value = 0.07403433476394849
$ws.Range("D8").Value = 0.01394849785407725
$ws.Range("E8").Value = 0.002145922746781116
$ws.Range("F8").Value = 0.06437768240343347
$ws.Range("J8").Value = 0.1223175965665236
$ws.Range("O8").Value = 0.01716738197424893
$ws.Range("Q8").Value = 0.167381974248927
$ws.Range("R8").Value = 0.09334763948497854
$ws.Range("S8").Value = 0.4452789699570815
$ws.Range("B9").Value = 0.09154929577464789
$ws.Range("D9").Value = 0.01056338028169014
$ws.Range("F9").Value = 0.06338028169014084
$ws.Range("J9").Value = 0.1056338028169014
$ws.Range("O9").Value = 0.0352112676056338
$ws.Range("Q9").Value = 0.1866197183098592
$ws.Range("R9").Value = 0.1232394366197183
$ws.Range("S9").Value = 0.3838028169014084
$ws.Range("B10").Value = 0.1039189784236019
$ws.Range("D10").Value = 0.01981505944517833
$ws.Range("E10").Value = 0.001761338617349185
$ws.Range("F10").Value = 0.07001321003963012
$ws.Range("J10").Value = 0.1096433289299868
$ws.Range("O10").Value = 0.01365037428445619
$ws.Range("Q10").Value = 0.2016732716864817
$ws.Range("R10").Value = 0.08718626155878467
$ws.Range("S10").Value = 0.3923381770145311
$ws.Range("G11").Value = 0.1882352941176471
$ws.Range("J11").Value = 0.08403361344537816
$ws.Range("K11").Value = 0.2403361344537815
$ws.Range("L11").Value = 0.4739495798319328
$ws.Range("S11").Value = 0.01344537815126051
$ws.Range("G12").Value = 0.7619047619047619
$ws.Range("J12").Value = 0.1802721088435374
$ws.Range("K12").Value = 0.01360544217687075
$ws.Range("L12").Value = 0.02040816326530612
$ws.Range("S12").Value = 0.02380952380952381
$ws.Range("G13").Value = 0.7333333333333333
$ws.Range("J13").Value = 0.2555555555555555
$ws.Range("S13").Value = 0.01111111111111111
$ws.Range("F15").Value = 0.01678657074340528
$ws.Range("H15").Value = 0.1990407673860911
$ws.Range("I15").Value = 0.05275779376498801
$ws.Range("J15").Value = 0.3597122302158273
$ws.Range("K15").Value = 0.09112709832134293
$ws.Range("M15").Value = 0.007194244604316547
$ws.Range("N15").Value = 0.002398081534772182
$ws.Range("O15").Value = 0.07434052757793765
$ws.Range("S15").Value = 0.1966426858513189
$ws.Range("F16").Value = 0.03115264797507788
$ws.Range("H16").Value = 0.1931464174454829
$ws.Range("I16").Value = 0.04361370716510903
$ws.Range("J16").Value = 0.3925233644859813
$ws.Range("K16").Value = 0.1246105919003115
$ws.Range("M16").Value = 0.02180685358255452
$ws.Range("N16").Value = 0.003115264797507788
$ws.Range("O16").Value = 0.02803738317757009
$ws.Range("S16").Value = 0.161993769470405
$ws.Range("F17").Value = 0.0262828535669587
$ws.Range("H17").Value = 0.2027534418022528
$ws.Range("I17").Value = 0.07509386733416772
$ws.Range("J17").Value = 0.3942428035043805
$ws.Range("K17").Value = 0.08760951188986232
$ws.Range("M17").Value = 0.03003754693366708
$ws.Range("O17").Value = 0.06257822277847309
$ws.Range("S17").Value = 0.1214017521902378
$ws.Range("F18").Value = 0.02412868632707775
$ws.Range("H18").Value = 0.2171581769436997
$ws.Range("I18").Value = 0.07506702412868632
$ws.Range("J18").Value = 0.3941018766756032
$ws.Range("K18").Value = 0.08042895442359249
$ws.Range("M18").Value = 0.01072386058981233
$ws.Range("O18").Value = 0.08310991957104558
$ws.Range("S18").Value = 0.1152815013404826
$ws.Range("F19").Value = 0.02237878159966846
$ws.Range("H19").Value = 0.226688769167012
$ws.Range("I19").Value = 0.06713634479900539
$ws.Range("J19").Value = 0.3622047244094488
$ws.Range("K19").Value = 0.1098217985909656
$ws.Range("M19").Value = 0.02320762536261915
$ws.Range("O19").Value = 0.06920845420638209
$ws.Range("S19").Value = 0.1193535018648985
